# Auto update Excel log
# Appends newly-logged sensor readings to the PIR, Humidity and Temperature sheets.

$wb = $excel.ActiveWorkbook

function Append-Rows {
    param($ws, $startRow, $rows, $textCols)
    for ($i = 0; $i -lt $rows.Length; $i++) {
        $r = $startRow + $i
        $rowData = $rows[$i]
        foreach ($tc in $textCols) {
            $ws.Cells.Item($r, $tc).NumberFormat = "@"
        }
        for ($j = 0; $j -lt $rowData.Length; $j++) {
            $cell = $ws.Cells.Item($r, $j + 1)
            $cell.Value = $rowData[$j]
        }
    }
}

$pirRows = @(
        @("2026-02-06", "10:05:57", "10:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-02-06", "10:06:00", "10:00", "Bathroom", "Motion Detected", "Active"),
        @("2026-02-06", "10:06:03", "10:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-02-06", "10:06:09", "10:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-02-06", "10:06:13", "10:00", "Bathroom", "Motion Detected", "Active"),
        @("2026-02-06", "10:06:21", "10:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-02-06", "10:06:22", "10:00", "Bathroom", "Motion Detected", "Active"),
        @("2026-02-06", "10:06:29", "10:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-02-06", "10:06:30", "10:00", "Bathroom", "Motion Detected", "Active"),
        @("2026-02-06", "10:06:37", "10:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-02-06", "10:06:38", "10:00", "Bathroom", "Motion Detected", "Active"),
        @("2026-02-06", "10:06:48", "10:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-02-06", "10:06:48", "10:00", "Bathroom", "Motion Detected", "Active"),
        @("2026-02-06", "10:06:54", "10:00", "Bathroom", "No Motion", "Inactive"),
        @("2026-02-06", "10:06:55", "10:00", "Bathroom", "Motion Detected", "Active")
)

$humidityRows = @(
        @("2026-02-06", "10:05:56", "10:00", "Bathroom", "70.0%", "Active"),
        @("2026-02-06", "10:05:58", "10:00", "Bathroom", "70.0%", "Active"),
        @("2026-02-06", "10:06:01", "10:00", "Bathroom", "70.0%", "Active"),
        @("2026-02-06", "10:06:06", "10:00", "Bathroom", "69.0%", "Active"),
        @("2026-02-06", "10:06:11", "10:00", "Bathroom", "70.0%", "Active"),
        @("2026-02-06", "10:06:16", "10:00", "Bathroom", "69.1%", "Active"),
        @("2026-02-06", "10:06:26", "10:00", "Bathroom", "70.1%", "Active"),
        @("2026-02-06", "10:06:36", "10:00", "Bathroom", "69.4%", "Active"),
        @("2026-02-06", "10:06:41", "10:00", "Bathroom", "69.4%", "Active"),
        @("2026-02-06", "10:06:46", "10:00", "Bathroom", "69.5%", "Active"),
        @("2026-02-06", "10:06:51", "10:00", "Bathroom", "69.4%", "Active"),
        @("2026-02-06", "10:06:56", "10:00", "Bathroom", "68.5%", "Active")
)

$temperatureRows = @(
        @("2026-02-06", "10:05:57", "10:00", "Bathroom", "27.8C", "Active"),
        @("2026-02-06", "10:05:59", "10:00", "Bathroom", "27.8C", "Active"),
        @("2026-02-06", "10:06:01", "10:00", "Bathroom", "27.8C", "Active"),
        @("2026-02-06", "10:06:06", "10:00", "Bathroom", "27.8C", "Active"),
        @("2026-02-06", "10:06:11", "10:00", "Bathroom", "27.8C", "Active"),
        @("2026-02-06", "10:06:17", "10:00", "Bathroom", "27.8C", "Active"),
        @("2026-02-06", "10:06:27", "10:00", "Bathroom", "27.9C", "Active"),
        @("2026-02-06", "10:06:37", "10:00", "Bathroom", "27.9C", "Active"),
        @("2026-02-06", "10:06:42", "10:00", "Bathroom", "27.9C", "Active"),
        @("2026-02-06", "10:06:47", "10:00", "Bathroom", "27.9C", "Active"),
        @("2026-02-06", "10:06:52", "10:00", "Bathroom", "27.9C", "Active")
)

$wsPir = $wb.Worksheets.Item("PIR")
Append-Rows $wsPir 311 $pirRows @(1)

$wsHumidity = $wb.Worksheets.Item("Humidity")
Append-Rows $wsHumidity 205 $humidityRows @(1, 5)

$wsTemperature = $wb.Worksheets.Item("Temperature")
Append-Rows $wsTemperature 205 $temperatureRows @(1)

Write-Host "Appended" $pirRows.Length "rows to PIR,"  $humidityRows.Length "rows to Humidity, and" $temperatureRows.Length "rows to Temperature."
